# Update product price cells (column E) on the "productos" sheet.
# Each edited cell's value changes and its style index shifts from the
# "plain currency" xf (14) to the "currency with explicit alignment" xf
# (13) -- the same pair of xfs Google Sheets toggles between when a cell's
# content is directly edited. Touching Orientation with its current/default
# value (0) is enough to make the host re-stamp the cell's alignment
# without altering how it actually looks, landing it on the matching xf.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("productos")

$updates = [ordered]@{
    "E70" = 56.05
    "E71" = 20.21
    "E72" = 36.24
    "E73" = 8.15
    "E74" = 14.6
    "E76" = 144.96
    "E77" = 55.63
    "E78" = 20.06
    "E79" = 35.97
    "E80" = 8.09
    "E81" = 14.39
    "E82" = 32.09
    "E83" = 143.87
    "E84" = 53.12
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.Value = $updates[$addr]
    $cell.Orientation = 0
}
